$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.080.24"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "3.326.84"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'581.34"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").Value = "'185.19"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.321.75"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "'47.05"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "'652.45"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("D15").Value = "3.855.50"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "'8.49"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "66.163.86"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'17.93"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "3.323.11"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "'17.91"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'100.83"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'5.03"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").Value = "'3.98"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "'9.49"
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("D29").Value = "'31.31"
$ws.Range("E29").Value = "  +2.87%  "
$ws.Range("D30").Value = "'8.46"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "'594.44"
$ws.Range("E32").Value = "  +3.25%  "
$ws.Range("D33").Value = "'3.86"
$ws.Range("E33").Value = "  -5.75%  "
$ws.Range("D34").Value = "'10.99"
$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "3.833.64"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'55.94"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'2.67"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0699"
$ws.Range("E40").Value = "  -4.94%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.126"
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("D42").Value = "'32.77"
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("E43").Value = "  +4.67%  "
$ws.Range("E44").Value = "  -5.34%  "
$ws.Range("D45").Value = "'0.335"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("D46").Value = "'0.0411"
$ws.Range("E46").Value = "  -3.93%  "
$ws.Range("D47").Value = "'3.04"
$ws.Range("E47").Value = "  -12.80%  "
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "'2.54"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").Value = "'130.51"
$ws.Range("E51").Value = "  +6.40%  "
